$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must remain text exactly as scraped (e.g. "1.00", "0.607",
# "3.470.13"). Excel auto-converts plain numeric-looking strings into Numbers, which
# both changes their cell type and can introduce floating point rounding artifacts
# (e.g. "0.607" -> 0.60699999999999998). Force the cell to Text format first, assign
# the literal string, then restore the default "Normal" style so no stray number-format
# styling is left behind on cells that originally had none.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.116.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.461.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.470.13"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.057.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.140.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.438.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0785"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.908.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0319"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.782"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "323.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.877"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("D51").Style = "Normal"

# Columns B (Coin), C (Link) and E (Volume) are never numeric-looking, so a plain
# value assignment is sufficient and keeps things simple.
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +5.20%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +3.64%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +10.12%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +5.82%  "
$ws.Range("E35").Value = "  +11.40%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +5.04%  "
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("E39").Value = "  +7.11%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E41").Value = "  +8.00%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("E46").Value = "  +11.80%  "
$ws.Range("E47").Value = "  +11.16%  "
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("E51").Value = "  +0.27%  "
